# Update review with PGx team:
# Mutant peak-height threshold for CYP2D6_001 (CYP2D6_14) raised from 300 to 1000,
# which means the detected mutant peak (height 428) no longer clears the
# threshold -> sample becomes wildtype instead of heterozygous for that marker,
# and the overall genotype result is recorded.

$wb = $excel.ActiveWorkbook

$peak = $wb.Worksheets.Item("peak_table")
$allele = $wb.Worksheets.Item("allele_table")
$marker = $wb.Worksheets.Item("marker_table")
$result = $wb.Worksheets.Item("genotype_result")

# peak_table: m_height for row 2 (CYP2D6_001 / CYP2D6_14) 300 -> 1000
$peak.Range("O2").Value = 1000

# allele_table: row 3 is the CYP2D6_001 / CYP2D6_14 mutant-base record
$allele.Range("K3").Value = 1000
$allele.Range("M3").Value = $false
$allele.Range("N3").Value = ""
$allele.Range("O3").Value = ""
$allele.Range("P3").Value = ""
$allele.Range("Q3").Value = ""
$allele.Range("R3").Value = "Peak(s) could not be detected. Please check peak ranges if required!"

# marker_table: row 2 genotype/phenotype for this marker
$marker.Range("G2").Value = "G"
$marker.Range("H2").Value = "wildtype"

# genotype_result: final diplotype for the sample
$result.Range("B2").Value = "*5/*10B"
